$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.366.82'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '3.185.79'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.37'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.27'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.182.31'
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.548'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.158'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.80'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -5.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.508'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000264'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.78'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').Value = '3.711.77'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').Value = '66.457.83'
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.39'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('D18').Value = '3.188.55'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '511.63'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.44'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.732'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.11'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.87'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.50'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.17'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.39'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +6.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.06'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +7.13%  '
$ws.Range('E31').Value = '  +5.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.07'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('B34').Value = 'Mantle'
$ws.Range('C34').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.21'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.54'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '510.37'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.73'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0894'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0423'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('E40').Value = '  +6.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.85'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.300'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.09%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.0₃0672'
$ws.Range('E43').Value = '  +5.12%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.83'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -5.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.43'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('D46').Value = '2.847.94'
$ws.Range('E46').Value = '  -5.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.34'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.42%  '
$ws.Range('E48').Value = '  +5.02%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.116'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.37%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.40'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.61%  '
